$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# A new daily price record (date 45050 = 2023-05-04, price 28.44) was
# added above the existing row for 45049, so insert a fresh row at 167
# and push the rest of the table down by one row.
$ws.Rows(167).Insert() | Out-Null

$ws.Range("A167").Value = 45050
$ws.Range("B167").Value = 28.44

# Reflect the cursor position that was active when the sheet was saved.
$ws.Range("E167").Select() | Out-Null
